$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.54"
$ws.Range("E2").Value = "'-0.99%"
$ws.Range("D3").Value = "'31.77"
$ws.Range("E3").Value = "'1.08%"
$ws.Range("D4").Value = "'5.091"
$ws.Range("E4").Value = "'-1.02%"
$ws.Range("D5").Value = "'0.08156"
$ws.Range("E5").Value = "'10.74%"
$ws.Range("D6").Value = "'2.523"
$ws.Range("E6").Value = "'-1.53%"
$ws.Range("D7").Value = "'7.769"
$ws.Range("E8").Value = "'2.40%"
$ws.Range("D9").Value = "'0.9303"
$ws.Range("E9").Value = "'1.03%"
$ws.Range("D10").Value = "'0.1759"
$ws.Range("E10").Value = "'0.82%"
$ws.Range("D11").Value = "'0.07513"
$ws.Range("E11").Value = "'0.27%"
$ws.Range("D12").Value = "'0.09002"
$ws.Range("E12").Value = "'10.59%"
$ws.Range("D13").Value = "'0.02997"
$ws.Range("E13").Value = "'-1.18%"
$ws.Range("E14").Value = "'0.87%"
$ws.Range("D15").Value = "'0.001506"
$ws.Range("E15").Value = "'0.50%"
$ws.Range("D16").Value = "'0.005760"
$ws.Range("E16").Value = "'-6.29%"
$ws.Range("D17").Value = "'3.582"
$ws.Range("E17").Value = "'3.86%"
$ws.Range("D18").Value = "'2.254"
$ws.Range("E18").Value = "'1.08%"
$ws.Range("E19").Value = "'-1.85%"
$ws.Range("E20").Value = "'-0.94%"
$ws.Range("D21").Value = "'3.925"
$ws.Range("E21").Value = "'-15.71%"
$ws.Range("D22").Value = "'0.1698"
$ws.Range("E22").Value = "'7.00%"
$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'-1.11%"
$ws.Range("E24").Value = "'1.16%"
$ws.Range("D25").Value = "'0.004458"
$ws.Range("E25").Value = "'-0.27%"
$ws.Range("E26").Value = "'-7.97%"
$ws.Range("D27").Value = "'0.0003409"
$ws.Range("E27").Value = "'82.14%"
$ws.Range("D39").Value = "'0.01764"
$ws.Range("E39").Value = "'1.94%"
$ws.Range("D40").Value = "'0.04526"
$ws.Range("E40").Value = "'-0.22%"
$ws.Range("D41").Value = "'0.006900"
$ws.Range("E41").Value = "'-4.30%"
$ws.Range("D42").Value = "'0.1354"
$ws.Range("E42").Value = "'0.68%"
$ws.Range("D43").Value = "'0.002207"
$ws.Range("E43").Value = "'-0.30%"
$ws.Range("D44").Value = "'0.009912"
$ws.Range("E44").Value = "'-9.17%"
$ws.Range("D45").Value = "'0.00006555"
$ws.Range("E45").Value = "'4.12%"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("E47").Value = "'-12.77%"
$ws.Range("E48").Value = "'11.08%"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("E50").Value = "'-0.06%"
